$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = "#06-파이썬(Python) 비교, 논리, 삼항연산자와 조건문"
$ws.Range("E4").Value = "https://teddylee777.github.io/python/python-tutorial-06"

$ws.Range("D5").Value = "비제차 미분방정식의 의미"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/25/nonhomogeneous_equation.html"

$ws.Range("D6").Value = "[Python] 날씨 시계열 데이터(Kaggle)로 ARIMA 적용하기"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-%EB%82%A0%EC%94%A8-%EC%8B%9C%EA%B3%84%EC%97%B4-%EB%8D%B0%EC%9D%B4%ED%84%B0Kaggle%EB%A1%9C-ARIMA-%EC%A0%81%EC%9A%A9%ED%95%98%EA%B8%B0"

$ws.Range("D9").Value = "[공지] 데이터 과학 대학원 입시 일정 + 6월 1일 2차 설명회"

$ws.Range("D16").Value = "Score-CAM : Score-weighted visual explanations for convolutional neural networks [XAI-9]"
$ws.Range("E16").Value = "https://wewinserv.tistory.com/151"

$ws.Range("D25").Value = "[바람돌이/딥러닝] GAN(4) - Improved Techniques for Training GANs 논문 이론 및 리뷰"
$ws.Range("E25").Value = "https://blog.naver.com/winddori2002/222365412213"

$ws.Range("D37").Value = "[Rehearsal] 2021 대한산업공학회 춘계 학술대회 - 김혜연"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1524&mod=document&pageid=1"

$ws.Range("D51").Value = "[python] python 환경 변수 설정하기 (AppData 폴더가 왜 없지 하시는 분들을 위해)"
$ws.Range("E51").Value = "https://bskyvision.com/1203"
